# Adds a new "2021" column (O) to the 15.b.1.1 environmental-protection
# table, mirroring the existing 2020 column (N) formatting for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year header): O4 = 2021, same look as N4 ---
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

# --- Row 5 (Total): O5 = 689 ---
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 689
$ws.Range("O5").NumberFormat = "0.0"
$ws.Range("O5").HorizontalAlignment = -4152

# --- Row 6: O6 = 94.1 ---
$ws.Range("N10").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 94.1

# --- Row 7: O7 = 147.1 ---
$ws.Range("N10").Copy()
$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("O7").Value = 147.1

# --- Row 8: O8 = 10.1 ---
$ws.Range("N10").Copy()
$ws.Range("O8").PasteSpecial(-4122)
$ws.Range("O8").Value = 10.1

# --- Row 9: O9 = "-" (text dash, shared with the rest of the row) ---
$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial(-4122)
$ws.Range("O9").Value = "-"
$ws.Range("O9").NumberFormat = "0.0"
$ws.Range("O9").HorizontalAlignment = -4152

# --- Row 10: O10 = 82.1 ---
$ws.Range("N10").Copy()
$ws.Range("O10").PasteSpecial(-4122)
$ws.Range("O10").Value = 82.1

# --- Row 11: O11 = 145.3 ---
$ws.Range("N10").Copy()
$ws.Range("O11").PasteSpecial(-4122)
$ws.Range("O11").Value = 145.30000000000001

# --- Row 12: O12 = 98.8 ---
$ws.Range("N10").Copy()
$ws.Range("O12").PasteSpecial(-4122)
$ws.Range("O12").Value = 98.8

# --- Row 13: O13 = 98.7 ---
$ws.Range("N10").Copy()
$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("O13").Value = 98.7

# --- Row 14: O14 = 1.8 ---
$ws.Range("N10").Copy()
$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("O14").Value = 1.8

# --- Row 15: O15 = "-" (text dash) ---
$ws.Range("N15").Copy()
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("O15").Value = "-"
$ws.Range("O15").NumberFormat = "0.0"
$ws.Range("O15").HorizontalAlignment = -4152

# --- Row 16: O16 = 10.9 (bottom-border row) ---
$ws.Range("N16").Copy()
$ws.Range("O16").PasteSpecial(-4122)
$ws.Range("O16").Value = 10.9
$ws.Range("O16").NumberFormat = "0.0"

# Match the saved cursor position recorded in the workbook.
[void]$ws.Range("P5").Select()
